$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("earnings_debt")

$values = @{
    "D" = 0.0434
    "E" = 0.161
    "G" = 0.09561817885649129
    "H" = 0.09561817885649129
    "I" = 0.1425313568985177
    "J" = 0.1142898496492758
    "K" = 51.9
    "L" = 0.08454145626323506
    "M" = 18.4
    "N" = 0.03857442348008385
    "O" = 0.3545279383429673
    "P" = 18.4
    "Q" = 0.03857442348008385
    "R" = 0.3545279383429673
    "U" = 592.1
    "V" = 1.241299790356394
    "W" = 0.1361847284177381
    "X" = 0.07573001268223215
    "Y" = 0.06045471573550598
    "Z" = -9.610507529979024
    "AA" = -1.098383460654536
    "AB" = 0.07572824868933122
    "AC" = -1.174111709343867
    "AD" = 0.015
    "AF" = 0.015
    "AG" = -592.085
    "AH" = 0.000031445552026665828386671475
    "AI" = 0.000034465723837643472463564132
    "AJ" = 5.144762566798452
    "AK" = 3.774006437836631
    "AN" = 0.0001706484641638225
    "AP" = -6.73589306029579
}

foreach ($col in $values.Keys) {
    $val = $values[$col]
    $ws.Range("$col`2").Value = $val
    $ws.Range("$col`3").Value = $val
}
